$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A33").Value = "1727-2097-1-ND"
$ws.Range("B33").Value = 1
$ws.Range("C33").Value = "D1"
$ws.Range("D2:D32").AutoFill($ws.Range("D2:D33"), 0)
$ws.Range("E33").Value = "3v3 zener diode "
